$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Conditional formatting: change the existing ">0.05" highlight rule
# into a 3-tier "<0.05 / <0.005" highlight, keeping the original rule
# (now "<0.05") at top priority and adding two new rules above it in
# the stack (">0.05" -> "<0.05", plus a new stronger "<0.005" rule).
$range = $ws.Range("B2:G31")

$original = $range.FormatConditions.Item(1)
$original.Operator = 6   # xlLess

$rule2 = $range.FormatConditions.Add(1, 6, "=0.05")
$rule2.Font.Color = -16383844
$rule2.Interior.Color = 13551615

$rule3 = $range.FormatConditions.Add(1, 6, "=0.005")
$rule3.Font.Color = -16383844
$rule3.Interior.Color = 13551615

# --- New note cell next to the table (single space placeholder) ---
$ws.Range("N16").Value = " "

# --- Restore the selection to the cell the author left active ---
$ws.Range("I13").Select()
